$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row: insert VAR(OK) before MSPE, rename Nugget -> S_nugget,
#     drop the old computed "Non-expl var of model" header and add two new
#     trailing headers VAR(TOTAL) and VAR(DATA) ---
$ws.Range("D1").Value = "VAR(OK)"
$ws.Range("E1").Value = "MSPE"
$ws.Range("F1").Value = "S_nugget"
$ws.Range("G1").Value = "VAR(TOTAL)"
$ws.Range("H1").Value = "VAR(DATA)"

# Give the two brand new header cells (G1,H1) the same bold/boxed style
# already used by the other header cells.
$ws.Range("F1").Copy() | Out-Null
$ws.Range("G1:H1").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# --- Data rows 2-13 ---
# columns: A, B, C, D(VAR_OK), E(MSPE), F(S_nugget), H(VAR_DATA)
# G(VAR_TOTAL) is left blank for every row.
$data = @(
    @(0,130000,1,0.3542437152146415,6.104421207852351,1.067,6.084087178003832),
    @(1,130000,2,0.5795627286992455,5.90785467806811,1.696,6.084087178003832),
    @(2,130000,3,0.7811871769957017,5.739350436129733,2.297,6.084087178003832),
    @(3,130000,4,0.9607850150456215,5.651751158741311,2.743,6.084087178003832),
    @(4,130000,5,1.0128557523289,5.641433115291756,3.063,6.084087178003832),
    @(5,130000,6,1.066439333051098,5.636778839762298,3.301,6.084087178003832),
    @(6,130000,7,1.082139319464439,5.634961651046167,3.508,6.084087178003832),
    @(7,130000,8,1.102006110258462,5.633633501109285,3.706,6.084087178003832),
    @(8,130000,9,1.123630734343489,5.62339082639934,3.813,6.084087178003832),
    @(9,130000,10,1.146000103602855,5.632906103503426,3.813,6.084087178003832),
    @(10,130000,11,1.154623899342249,5.63576651123089,3.813,6.084087178003832),
    @(11,130000,12,1.156557324888634,5.63690469376062,3.813,6.084087178003832)
)

$row = 2
foreach ($r in $data) {
    $ws.Cells.Item($row, 1).Value = $r[0]
    $ws.Cells.Item($row, 2).Value = $r[1]
    $ws.Cells.Item($row, 3).Value = $r[2]
    $ws.Cells.Item($row, 4).Value = $r[3]
    $ws.Cells.Item($row, 5).Value = $r[4]
    $ws.Cells.Item($row, 6).Value = $r[5]
    $ws.Cells.Item($row, 7).Value = ""
    $ws.Cells.Item($row, 8).Value = $r[6]
    $row = $row + 1
}

# Copy the style of the existing "#PCs label" column A cells down onto the
# three newly added rows (11-13) so they match rows 2-10.
$ws.Range("A10").Copy() | Out-Null
$ws.Range("A11:A13").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0
